$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 331; existing rows 331-349 shift down to 332-350.
$ws.Rows.Item(331).Insert()

# Populate the newly inserted row 331 with this week's data.
$ws.Range("A331").Value = 8
$ws.Range("B331").Value = "Terminal La Palmera de La Serena"
$ws.Range("C331").Value = "Coquimbo"
$ws.Range("D331").Value = 44706
$ws.Range("E331").Value = 4
$ws.Range("F331").Value = 100114013
$ws.Range("G331").Value = "Zanahoria"
$ws.Range("H331").Value = "Sin especificar"
$ws.Range("I331").Value = "Primera"
$ws.Range("J331").Value = 680
$ws.Range("K331").Value = 6000
$ws.Range("L331").Value = 7000
$ws.Range("M331").Value = 6500
$ws.Range("N331").Value = "`$/saco 20 kilos"
$ws.Range("O331").Value = "Provincia del Elquí"
$ws.Range("P331").Value = 325
$ws.Range("Q331").Value = 20
$ws.Range("R331").Value = "Hortaliza"

# Ensure the date cell keeps the same date/time number format used by the rest of column D.
$ws.Range("D331").NumberFormat = $ws.Range("D332").NumberFormat
